# "arreglo de archivo patrocinadores"
#
# The dropdown source list lives on the (previously hidden) Hoja2. The
# author added a new entry, "Patrocinadores", right before "Otros", then
# unhid the helper sheet again.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Unhide the helper sheet (Hoja2) that holds the dropdown source list.
$ws2.Visible = -1   # xlSheetVisible

# Insert a new row for "Patrocinadores" right above "Otros" (row 6),
# pushing "Otros" down to row 7.
$ws2.Rows.Item(6).Insert()
$ws2.Range("A6").Value = "Patrocinadores"

# Column A is now wider to fit the longest entry.
$ws2.Columns.Item(1).AutoFit()

# Restore the cursor positions seen in the saved file: Hoja2 selection on
# C5, then back to Hoja1 (which stays the active/tab-selected sheet) at F7.
$ws2.Range("C5").Select()
$ws1.Range("F7").Select()
